# edit.ps1 - apply the ldmc.docx revision:
#   * Model formula text: sqrt(LDMC) -> (LDMC)^(1/3)            (3 paragraphs)
#   * Mojibake chi character fix: "Ï‡" -> "χ"                    (5 table headers)
#   * A handful of updated numeric results in the 5 result tables
#   * A few column-width / row-height tweaks on the tables
#
# Helper: replace the text inside a single table cell without disturbing
# any other cell containing the same literal text elsewhere in the
# document. (Scoping Find.Execute to $cell.Range directly does not
# commit the edit in this host, but re-wrapping the same start/end
# offsets in a document Range does.)
function Set-CellText($cell, $oldText, $newText) {
    $cr = $cell.Range
    $rng = $d.Range($cr.Start, $cr.End)
    # Replace:=1 (wdReplaceOne) -- wdReplaceAll loops forever in this host
    # whenever the replacement text itself contains the search text
    # (e.g. "0" -> "0.037"), so we only ever replace the single match
    # that the (tightly cell-scoped) range can contain anyway.
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $newText, 1) | Out-Null
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Model: sqrt(LDMC) ~ ..." -> "Model: (LDMC)^(1/3) ~ ..." (all 3 spots)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("sqrt(LDMC)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(LDMC)^(1/3)", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Mojibake "Ï‡" -> proper "χ" in every table header (5 occurrences)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Ï‡", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "χ", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Table 1 ("Table 1: Test for variance among families and populations")
# ---------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$t1.Rows.Item(1).Height = 615 / 20.0
$t1.Cell(1, 6).Width = 1084 / 20.0

Set-CellText $t1.Cell(3, 3) "0.003" "0.093"
Set-CellText $t1.Cell(3, 5) "0.066" "0.407"
Set-CellText $t1.Cell(3, 6) "0.479" "0.3805"
Set-CellText $t1.Cell(4, 4) "0.003" "0.002"
Set-CellText $t1.Cell(4, 5) "99.934" "99.593"

# ---------------------------------------------------------------------
# 4) Table 2 ("Table 2: Assess how much variance is explained by urbanization")
# ---------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$t2.Rows.Item(1).Height = 615 / 20.0
$t2.Cell(1, 3).Width = 961 / 20.0
$t2.Cell(1, 5).Width = 1084 / 20.0
$t2.Cell(1, 6).Width = 961 / 20.0

Set-CellText $t2.Cell(2, 3) "0" "0.000"
Set-CellText $t2.Cell(2, 5) "0" "0.000"
Set-CellText $t2.Cell(3, 3) "0" "0.037"
Set-CellText $t2.Cell(3, 5) "0" "0.254"
Set-CellText $t2.Cell(3, 6) "0.5" "0.424"
Set-CellText $t2.Cell(4, 4) "0.003" "0.002"
Set-CellText $t2.Cell(4, 5) "100" "99.746"

# ---------------------------------------------------------------------
# 5) Table 3 ("Table 3: Quantify variance explained by urbanization")
# ---------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$t3.Rows.Item(1).Height = 571 / 20.0

Set-CellText $t3.Cell(2, 3) "47.653" "58.087"
Set-CellText $t3.Cell(3, 3) "2.292" "2.091"
Set-CellText $t3.Cell(3, 4) "0.13" "0.148"

# ---------------------------------------------------------------------
# 6) Table 4 ("Table 4: Assess how much variance is explained by urbanization")
# ---------------------------------------------------------------------
$t4 = $d.Tables.Item(4)
$t4.Rows.Item(1).Height = 615 / 20.0
$t4.Cell(1, 6).Width = 961 / 20.0

Set-CellText $t4.Cell(3, 3) "0.012" "0.157"
Set-CellText $t4.Cell(3, 5) "0.138" "0.537"
Set-CellText $t4.Cell(3, 6) "0.4565" "0.346"
Set-CellText $t4.Cell(4, 4) "0.003" "0.002"
Set-CellText $t4.Cell(4, 5) "99.862" "99.463"

# ---------------------------------------------------------------------
# 7) Table 5 ("Table 5: Quantify variance explained by urbanization")
# ---------------------------------------------------------------------
$t5 = $d.Tables.Item(5)
$t5.Rows.Item(1).Height = 571 / 20.0

Set-CellText $t5.Cell(2, 3) "47.636" "58.118"
Set-CellText $t5.Cell(3, 3) "0.355" "0.140"
Set-CellText $t5.Cell(3, 4) "0.551" "0.708"
